$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Role column (F) values -------------------------------------------------
# Fill in the previously-empty Role cells for the three testers first so the
# shared-string table grows in the same order as the target file, then give
# the former "Leader" cell (F2) its new, longer role text.
$ws.Range("F3").Value = "Tester, focusing on UI automation and functional testing"
$ws.Range("F4").Value = "Tester, specializing in API testing and integration validation"
$ws.Range("F5").Value = "Tester, handling performance testing and optimization"

# F2 switches from the bold/white header-style look to the same plain
# "name/role" style used by column C (style index 4) before getting its new
# text - copy C2's formatting across, then overwrite the value.
$ws.Range("C2").Copy()
$ws.Range("F2").PasteSpecial(-4122)
$excel.CutCopyMode = $False
$ws.Range("F2").Value = "Scrum Master & Team Lead, responsible for coordination and removing blockers"

# --- Sheet view -------------------------------------------------------------
# Turn off right-to-left sheet display, zoom in to 130%, and move the
# selection from D11 to C6.
$excel.ActiveWindow.DisplayRightToLeft = $False
$excel.ActiveWindow.Zoom = 130
$ws.Range("C6").Select()

# --- Column width / row height ----------------------------------------------
# Column F widens considerably now that it holds long role descriptions.
$ws.Columns("F").ColumnWidth = 69.04
# Row 2's height shrinks slightly to match the other data rows.
$ws.Rows(2).RowHeight = 15
